$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 380, shifting existing rows 380:447 down to 381:448
$ws.Rows.Item(380).Insert()

# Populate the new row 380 with the new data entry
$ws.Cells.Item(380, 1).Value = 10
$ws.Cells.Item(380, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(380, 3).Value = "La Araucanía"
$ws.Cells.Item(380, 4).Value = 45015
$ws.Cells.Item(380, 5).Value = 9
$ws.Cells.Item(380, 6).Value = 100114013
$ws.Cells.Item(380, 7).Value = "Zanahoria"
$ws.Cells.Item(380, 8).Value = "Sin especificar"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 350
$ws.Cells.Item(380, 11).Value = 6000
$ws.Cells.Item(380, 12).Value = 7000
$ws.Cells.Item(380, 13).Value = 6571
$ws.Cells.Item(380, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(380, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(380, 16).Value = 263
$ws.Cells.Item(380, 17).Value = 25
$ws.Cells.Item(380, 18).Value = "Hortaliza"
